$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy date/time formatting from the last existing row (A113) onto the new date cells
$ws.Range("A113").Copy()
$ws.Range("A114:A115").PasteSpecial(-4122)

# Row 114
$ws.Cells.Item(114, 1).Value = 45477.2916666667
$ws.Cells.Item(114, 2).Value = 0
$ws.Cells.Item(114, 3).Value = 1.91999995708466
$ws.Cells.Item(114, 4).Value = 1.91999995708466
$ws.Cells.Item(114, 5).Value = 1.91999995708466
$ws.Cells.Item(114, 6).Value = 1.91999995708466
$ws.Cells.Item(114, 7).NumberFormat = "@"
$ws.Cells.Item(114, 7).Value = "1.91999995708466"
$ws.Cells.Item(114, 7).Style = "Normal"
$ws.Cells.Item(114, 8).Value = "KK.MI"

# Row 115
$ws.Cells.Item(115, 1).Value = 45478.5668055556
$ws.Cells.Item(115, 2).Value = 900
$ws.Cells.Item(115, 3).Value = 1.87000000476837
$ws.Cells.Item(115, 4).Value = 1.87000000476837
$ws.Cells.Item(115, 5).Value = 1.87000000476837
$ws.Cells.Item(115, 6).Value = 1.87000000476837
$ws.Cells.Item(115, 7).NumberFormat = "@"
$ws.Cells.Item(115, 7).Value = "1.87000000476837"
$ws.Cells.Item(115, 7).Style = "Normal"
$ws.Cells.Item(115, 8).Value = "KK.MI"
